$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61 (shifts existing rows 61..91 down to 62..92,
# carrying their data/formatting along automatically).
$ws.Rows(61).Insert()

# Populate the newly inserted row 61 with the new weekly price record.
$ws.Range("A61").Value = 5
$ws.Range("B61").Value = "Macroferia Regional de Talca"
$ws.Range("C61").Value = "Maule"
$ws.Range("D61").Value = 44609
$ws.Range("E61").Value = 7
$ws.Range("F61").Value = 100112022
$ws.Range("G61").Value = "Arveja Verde"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 300
$ws.Range("K61").Value = 28000
$ws.Range("L61").Value = 28000
$ws.Range("M61").Value = 28000
$ws.Range("N61").Value = "$/saco 25 kilos"
$ws.Range("O61").Value = "Carahue"
$ws.Range("P61").Value = 1120
$ws.Range("Q61").Value = 25
$ws.Range("R61").Value = "Hortaliza"
